$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly scoreboard rows (Week 6) appended below the existing data (rows 220-230).
# Columns: A Participant, B Date, C Workout Type, D Total Duration, E Total Distance,
#          F Total Elevation, G Zone 1, H Zone 2, I Zone 3, J Zone 4, K Zone 5,
#          L Workout Level, M Week

$rows = @(
    @("Eric",     45489, "Workout", 55, 0,     0,   23, 30, 2,  0, 0,  "Wily Hyena",       6),
    @("Steven",   45489, "Workout", 38, 0,     0,   8,  18, 12, 1, 0,  "Brave Leopard",    6),
    @("Steven",   45489, "Walk",    30, 1.44,  98,  30, 0,  0,  0, 0,  "Brave Leopard",    6),
    @("Jeremiah", 45490, "Workout", 41, 0,     0,   40, 1,  0,  0, 0,  "Wily Hyena",       6),
    @("Eric",     45490, "Workout", 30, 0,     0,   21, 10, 0,  0, 0,  "Wily Hyena",       6),
    @("Phil",     45490, "Workout", 22, 0,     0,   22, 0,  0,  0, 0,  "Sauntering Hippo", 6),
    @("Matt",     45490, "Ride",    63, 20.01, 0,   18, 46, 0,  0, 0,  "Wily Hyena",       6),
    @("Steven",   45491, "Run",     26, 2.39,  135, 1,  6,  15, 4, 0,  "Brave Leopard",    6),
    @("Steven",   45491, "Walk",    28, 1.15,  85,  27, 1,  0,  0, 0,  "Brave Leopard",    6),
    @("Matt",     45491, "Run",     48, 5.41,  361, 0,  2,  10, 8, 25, "Wily Hyena",       6),
    @("Matt",     45491, "Walk",    10, 0.3,   23,  10, 0,  0,  0, 0,  "Wily Hyena",       $null)
)

$startRow = 220
$lastRow  = $startRow + $rows.Count - 1

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value  = $data[0]           # A Participant
    $ws.Cells.Item($r, 2).Value  = $data[1]           # B Date
    $ws.Cells.Item($r, 3).Value  = $data[2]           # C Workout Type
    $ws.Cells.Item($r, 4).Value  = $data[3]           # D Total Duration
    $ws.Cells.Item($r, 5).Value  = $data[4]           # E Total Distance
    $ws.Cells.Item($r, 6).Value  = $data[5]           # F Total Elevation
    $ws.Cells.Item($r, 7).Value  = $data[6]           # G Zone 1
    $ws.Cells.Item($r, 8).Value  = $data[7]           # H Zone 2
    $ws.Cells.Item($r, 9).Value  = $data[8]           # I Zone 3
    $ws.Cells.Item($r, 10).Value = $data[9]           # J Zone 4
    $ws.Cells.Item($r, 11).Value = $data[10]          # K Zone 5
    $ws.Cells.Item($r, 12).Value = $data[11]          # L Workout Level

    if ($null -ne $data[12]) {
        $ws.Cells.Item($r, 13).Value = $data[12]      # M Week
    }
}

# Apply the existing date-column number format (style index 1, numFmt 14) to
# the new Date cells by copying formatting from the last pre-existing date
# cell, matching the source workbook instead of inventing a new style.
$ws.Cells.Item(219, 2).Copy()
$ws.Range($ws.Cells.Item($startRow, 2), $ws.Cells.Item($lastRow, 2)).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the sheet view to match: scroll the frozen pane down and move the
# active selection onto the newly entered data.
$ws.Application.ActiveWindow.ScrollRow = 205
$ws.Range("L221").Select()
